# Animation Retargeting Tutorial Renewal - PPT data update
#
# Targets slide 4 ("Animator Controller Example") of the deck:
#   - Shrink the width of the text box that sits under the "Base Layer"
#     screenshot (the animation-speed-button paragraph box) slightly,
#     and pin down an explicit 0 rotation on it.
#   - Drop the leftover local-machine `descr` (alt text) that Polaris
#     Office had stamped onto two re-embedded screenshots, now that the
#     underlying image resources were deleted/re-added.
#
# (The companion bump of the <p:sldMasterId> internal id and the
# r:id renumbering of three <p:sldId> entries in ppt/presentation.xml
# are bookkeeping byproducts of PowerPoint's own package serializer
# from deleting/re-adding resources elsewhere in the authoring
# session - there is no Slides/Shapes object-model property that
# addresses a relationship id or the master's internal id directly,
# so nothing in this script targets them.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# "텍스트 상자 24" (id 1177) - resize width 324.95pt -> 324.25pt (4126865 -> 4117975 EMU)
# and stamp an explicit Rotation of 0 (writes rot="0" on the <a:xfrm>).
$textBox = $s.Shapes.Item(3)
$textBox.Width = 324.25
$textBox.Rotation = 0

# "그림 1" (id 1180) - clear the stale local-path alt text / descr.
$pic1 = $s.Shapes.Item(6)
$pic1.AlternativeText = ""

# "그림 5" (id 1182) - clear the stale local-path alt text / descr.
$pic5 = $s.Shapes.Item(8)
$pic5.AlternativeText = ""
